# Applies the scheduled-runner value refresh to the Leviathan_Profits workbook.
# For each affected Leve row, the market-price-derived columns (H:N) are recomputed;
# a handful of rows also gain/lose their LeveProfitHQ (N) cell depending on whether an
# HQ price is available for that leve, mirroring the upstream source-data structure.
$wb = $excel.ActiveWorkbook

# ===== Sheet: ALC =====
$ws = $wb.Worksheets.Item("ALC")
# Row 6
$ws.Range("H6").Value = 4443.4546
$ws.Range("I6").Value = 348.16666
$ws.Range("J6").Value = 9357.799999999999
$ws.Range("K6").Value = 1044.49998
$ws.Range("L6").Value = 28073.4
$ws.Range("M6").Value = -932.4999800000001
$ws.Range("N6").Value = -28297.4

# Row 19
$ws.Range("H19").Value = 5029.1665
$ws.Range("I19").Value = 6095
$ws.Range("J19").Value = 3963.3333
$ws.Range("K19").Value = 6095
$ws.Range("L19").Value = 3963.3333
$ws.Range("M19").Value = -5920
$ws.Range("N19").Value = -4313.3333

# Row 100
$ws.Range("H100").Value = 3275
$ws.Range("I100").Value = 3275
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 3275
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -2734
$ws.Range("N100").ClearContents()

# Row 106
$ws.Range("H106").Value = 8957.333000000001
$ws.Range("I106").Value = 1838.9
$ws.Range("J106").Value = 23194.2
$ws.Range("K106").Value = 1838.9
$ws.Range("L106").Value = 23194.2
$ws.Range("M106").Value = -1207.9
$ws.Range("N106").Value = -24456.2

# Row 132
$ws.Range("H132").Value = 3146.617
$ws.Range("I132").Value = 1303.4054
$ws.Range("J132").Value = 9966.5
$ws.Range("K132").Value = 3910.2162
$ws.Range("L132").Value = 29899.5
$ws.Range("M132").Value = -1380.2162

# Row 138
$ws.Range("H138").Value = 1747.8214
$ws.Range("I138").Value = 1073.2632
$ws.Range("J138").Value = 3171.889
$ws.Range("K138").Value = 3219.7896
$ws.Range("L138").Value = 9515.667000000001
$ws.Range("M138").Value = 1920.2104
$ws.Range("N138").Value = -19795.667

# ===== Sheet: ARM =====
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 36457.93
$ws.Range("I32").Value = 21583.959
$ws.Range("J32").Value = 117438.445
$ws.Range("K32").Value = 21583.959
$ws.Range("L32").Value = 117438.445
$ws.Range("M32").Value = -21296.959
$ws.Range("N32").Value = -118012.445

# Row 45
$ws.Range("H45").Value = 14235.814
$ws.Range("I45").Value = 13235.579
$ws.Range("J45").Value = 16611.375
$ws.Range("K45").Value = 13235.579
$ws.Range("L45").Value = 16611.375
$ws.Range("M45").Value = -12858.579

# Row 61
$ws.Range("H61").Value = 4573.3
$ws.Range("I61").Value = 4497.3335
$ws.Range("J61").Value = 4687.25
$ws.Range("K61").Value = 4497.3335
$ws.Range("L61").Value = 4687.25
$ws.Range("M61").Value = -4285.3335
$ws.Range("N61").Value = -5111.25

# Row 132
$ws.Range("H132").Value = 2081.3235
$ws.Range("I132").Value = 1582.65
$ws.Range("J132").Value = 2793.7144
$ws.Range("K132").Value = 4747.950000000001
$ws.Range("L132").Value = 8381.143199999999
$ws.Range("M132").Value = -2217.950000000001
$ws.Range("N132").Value = -13441.1432

# Row 136
$ws.Range("H136").Value = 4573.3
$ws.Range("I136").Value = 4497.3335
$ws.Range("J136").Value = 4687.25
$ws.Range("K136").Value = 13492.0005
$ws.Range("L136").Value = 14061.75
$ws.Range("M136").Value = -10942.0005
$ws.Range("N136").Value = -19161.75

# ===== Sheet: BSM =====
$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 6543986.5
$ws.Range("I20").Value = 13894086
$ws.Range("J20").Value = 10564.777
$ws.Range("K20").Value = 13894086
$ws.Range("L20").Value = 10564.777
$ws.Range("M20").Value = -13893839
$ws.Range("N20").Value = -11058.777

# Row 134
$ws.Range("H134").Value = 2048.3635
$ws.Range("I134").Value = 1942.25
$ws.Range("J134").Value = 3109.5
$ws.Range("K134").Value = 5826.75
$ws.Range("L134").Value = 9328.5
$ws.Range("M134").Value = -3291.75
$ws.Range("N134").Value = -14398.5

# ===== Sheet: CRP =====
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 4064
$ws.Range("I31").Value = 1513.8572
$ws.Range("J31").Value = 6047.4443
$ws.Range("K31").Value = 1513.8572
$ws.Range("L31").Value = 6047.4443
$ws.Range("M31").Value = -1218.8572

# Row 34
$ws.Range("H34").Value = 4064
$ws.Range("I34").Value = 1513.8572
$ws.Range("J34").Value = 6047.4443
$ws.Range("K34").Value = 1513.8572
$ws.Range("L34").Value = 6047.4443
$ws.Range("M34").Value = -1311.8572

# Row 59
$ws.Range("H59").Value = 16249.75
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 16249.75
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 16249.75
$ws.Range("N59").Value = -18539.75

# ===== Sheet: CUL =====
$ws = $wb.Worksheets.Item("CUL")
# Row 68
$ws.Range("H68").Value = 1685.25
$ws.Range("I68").Value = 1600
$ws.Range("J68").Value = 1713.6666
$ws.Range("K68").Value = 4800
$ws.Range("L68").Value = 5140.9998
$ws.Range("M68").Value = -3989
$ws.Range("N68").Value = -6762.9998

# Row 70
$ws.Range("H70").Value = 1677
$ws.Range("I70").Value = 1677
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 5031
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = -4716

# Row 71
$ws.Range("H71").Value = 1685.25
$ws.Range("I71").Value = 1600
$ws.Range("J71").Value = 1713.6666
$ws.Range("K71").Value = 14400
$ws.Range("L71").Value = 15422.9994
$ws.Range("M71").Value = -10344
$ws.Range("N71").Value = -23534.9994

# Row 73
$ws.Range("H73").Value = 1677
$ws.Range("I73").Value = 1677
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 5031
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = -3939

# Row 113
$ws.Range("H113").Value = 618.2
$ws.Range("I113").Value = 618.2
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1854.6
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 315.3999999999999
$ws.Range("N113").ClearContents()

# Row 131
$ws.Range("H131").Value = 1724.6154
$ws.Range("I131").Value = 1129.3334
$ws.Range("J131").Value = 2536.3635
$ws.Range("K131").Value = 3388.0002
$ws.Range("L131").Value = 7609.0905
$ws.Range("M131").Value = 1651.9998
$ws.Range("N131").Value = -17689.0905

# ===== Sheet: GSM =====
$ws = $wb.Worksheets.Item("GSM")
# Row 35
$ws.Range("H35").Value = 17411
$ws.Range("I35").Value = 17875
$ws.Range("J35").Value = 15555
$ws.Range("K35").Value = 17875
$ws.Range("L35").Value = 15555
$ws.Range("M35").Value = -17577
$ws.Range("N35").Value = -16151

# Row 132
$ws.Range("H132").Value = 1539.7142
$ws.Range("I132").Value = 1603.6666
$ws.Range("J132").Value = 1491.75
$ws.Range("K132").Value = 4810.9998
$ws.Range("L132").Value = 4475.25
$ws.Range("M132").Value = -2280.9998
$ws.Range("N132").Value = -9535.25

# ===== Sheet: LTW =====
$ws = $wb.Worksheets.Item("LTW")
# Row 93
$ws.Range("H93").Value = 1879.25
$ws.Range("I93").Value = 1916.8462
$ws.Range("J93").Value = 1716.3334
$ws.Range("K93").Value = 1916.8462
$ws.Range("L93").Value = 1716.3334
$ws.Range("M93").Value = -668.8462
$ws.Range("N93").Value = -4212.3334

# Row 132
$ws.Range("H132").Value = 2620.5667
$ws.Range("I132").Value = 2418.7827
$ws.Range("J132").Value = 3283.5715
$ws.Range("K132").Value = 7256.348100000001
$ws.Range("L132").Value = 9850.7145
$ws.Range("M132").Value = -4726.348100000001
$ws.Range("N132").Value = -14910.7145

# Row 133
$ws.Range("H133").Value = 79612.5
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 79612.5
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 79612.5
$ws.Range("N133").Value = -84672.5

# Row 136
$ws.Range("H136").Value = 3333.8235
$ws.Range("I136").Value = 2483.8572
$ws.Range("J136").Value = 3928.8
$ws.Range("K136").Value = 7451.571599999999
$ws.Range("L136").Value = 11786.4
$ws.Range("M136").Value = -4901.571599999999

# ===== Sheet: WVR =====
$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 1223.7142
$ws.Range("I122").Value = 1223.7142
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 3671.1426
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -1221.1426

# Row 132
$ws.Range("H132").Value = 10071
$ws.Range("I132").Value = 10197
$ws.Range("J132").Value = 9000
$ws.Range("K132").Value = 30591
$ws.Range("L132").Value = 27000
$ws.Range("M132").Value = -28061

# Row 135
$ws.Range("H135").Value = 0
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").Value = 0
